$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 796.3570999999999
$ws.Range("I6").Value = 109.47619
$ws.Range("J6").Value = 2857
$ws.Range("K6").Value = 328.42857
$ws.Range("L6").Value = 8571
$ws.Range("M6").Value = -216.42857
$ws.Range("N6").Value = -8795
$ws.Range("H41").Value = 105.5
$ws.Range("J41").Value = 116
$ws.Range("L41").Value = 116
$ws.Range("N41").Value = -996
$ws.Range("H52").Value = 3223.8462
$ws.Range("J52").Value = 3446.3635
$ws.Range("L52").Value = 10339.0905
$ws.Range("N52").Value = -10659.0905
$ws.Range("H76").Value = 3355.5557
$ws.Range("I76").Value = 3166.6667
$ws.Range("J76").Value = 3733.3333
$ws.Range("K76").Value = 3166.6667
$ws.Range("L76").Value = 3733.3333
$ws.Range("M76").Value = -2851.6667
$ws.Range("N76").Value = -4363.3333
$ws.Range("H79").Value = 3355.5557
$ws.Range("I79").Value = 3166.6667
$ws.Range("J79").Value = 3733.3333
$ws.Range("K79").Value = 3166.6667
$ws.Range("L79").Value = 3733.3333
$ws.Range("M79").Value = -2074.6667
$ws.Range("N79").Value = -5917.3333
$ws.Range("H135").Value = 1954.8422
$ws.Range("I135").Value = 1980.1111
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 17820.9999
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -15285.9999
$ws.Range("N135").Value = -18570
$ws.Range("H137").Value = 2933.963
$ws.Range("I137").Value = 2933.963
$ws.Range("K137").Value = 8801.889000000001
$ws.Range("M137").Value = -6251.889000000001
$ws.Range("H138").Value = 155429.22
$ws.Range("J138").Value = 189561.4
$ws.Range("L138").Value = 568684.2
$ws.Range("N138").Value = -578964.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 788334.5
$ws.Range("I32").Value = 1081766.8
$ws.Range("J32").Value = 14740.182
$ws.Range("K32").Value = 1081766.8
$ws.Range("L32").Value = 14740.182
$ws.Range("M32").Value = -1081479.8
$ws.Range("N32").Value = -15314.182
$ws.Range("H122").Value = 2257.2
$ws.Range("I122").Value = 1851.6
$ws.Range("J122").Value = 2662.8
$ws.Range("K122").Value = 5554.799999999999
$ws.Range("L122").Value = 7988.400000000001
$ws.Range("M122").Value = -3104.799999999999
$ws.Range("N122").Value = -12888.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 47620940
$ws.Range("I86").Value = 66668452
$ws.Range("J86").Value = 2166.6667
$ws.Range("K86").Value = 66668452
$ws.Range("L86").Value = 2166.6667
$ws.Range("M86").Value = -66667329
$ws.Range("N86").Value = -4412.6667
$ws.Range("H89").Value = 47620940
$ws.Range("I89").Value = 66668452
$ws.Range("J89").Value = 2166.6667
$ws.Range("K89").Value = 333342260
$ws.Range("L89").Value = 10833.3335
$ws.Range("M89").Value = -333336644
$ws.Range("N89").Value = -22065.3335

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 14777.333
$ws.Range("J15").Value = 19666
$ws.Range("L15").Value = 19666
$ws.Range("N15").Value = -20006
$ws.Range("H31").Value = 2213.6445
$ws.Range("I31").Value = 1046.1538
$ws.Range("K31").Value = 1046.1538
$ws.Range("M31").Value = -751.1538
$ws.Range("H34").Value = 2213.6445
$ws.Range("I34").Value = 1046.1538
$ws.Range("K34").Value = 1046.1538
$ws.Range("M34").Value = -844.1538
$ws.Range("H134").Value = 1688.8462
$ws.Range("I134").Value = 1574.3158
$ws.Range("J134").Value = 1999.7142
$ws.Range("K134").Value = 4722.9474
$ws.Range("L134").Value = 5999.142599999999
$ws.Range("M134").Value = -2187.9474
$ws.Range("N134").Value = -11069.1426

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 107
$ws.Range("I2").Value = 17.5
$ws.Range("J2").Value = 166.66667
$ws.Range("K2").Value = 105
$ws.Range("L2").Value = 1000.00002
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = -1226.00002
$ws.Range("H39").Value = 1755.5714
$ws.Range("J39").Value = 1999.8334
$ws.Range("L39").Value = 5999.5002
$ws.Range("N39").Value = -6587.5002
$ws.Range("H40").Value = 114.5625
$ws.Range("I40").Value = 110.5
$ws.Range("K40").Value = 442
$ws.Range("M40").Value = -373
$ws.Range("H68").Value = 1884
$ws.Range("I68").Value = 1884
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5652
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -4841
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1884
$ws.Range("I71").Value = 1884
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 16956
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -12900
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 1512
$ws.Range("I113").Value = 766.6667
$ws.Range("J113").Value = 2257.3333
$ws.Range("K113").Value = 2300.0001
$ws.Range("L113").Value = 6771.999899999999
$ws.Range("M113").Value = -130.0001000000002
$ws.Range("N113").Value = -11111.9999
$ws.Range("H122").Value = 7211.067
$ws.Range("I122").Value = 418.77777
$ws.Range("J122").Value = 17399.5
$ws.Range("K122").Value = 3768.99993
$ws.Range("L122").Value = 156595.5
$ws.Range("M122").Value = -1318.99993
$ws.Range("N122").Value = -161495.5
$ws.Range("H132").Value = 3106.8765
$ws.Range("I132").Value = 2127.75
$ws.Range("K132").Value = 19149.75
$ws.Range("M132").Value = -16619.75
$ws.Range("H136").Value = 2275.7144
$ws.Range("I136").Value = 1826
$ws.Range("K136").Value = 5478
$ws.Range("M136").Value = -378
$ws.Range("H139").Value = 4155.8
$ws.Range("I139").Value = 1284.2222
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 3852.6666
$ws.Range("L139").Value = 90000
$ws.Range("M139").Value = 1287.3334
$ws.Range("N139").Value = -100280
$ws.Range("H140").Value = 1658.3334
$ws.Range("I140").Value = 1336.5385
$ws.Range("J140").Value = 3750
$ws.Range("K140").Value = 4009.6155
$ws.Range("L140").Value = 11250
$ws.Range("M140").Value = 1170.3845
$ws.Range("N140").Value = -21610
$ws.Range("H141").Value = 3444
$ws.Range("I141").Value = 1082.1538
$ws.Range("J141").Value = 7282
$ws.Range("K141").Value = 3246.4614
$ws.Range("L141").Value = 21846
$ws.Range("M141").Value = 1933.5386
$ws.Range("N141").Value = -32206

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 428181.6
$ws.Range("I11").Value = 666966.7
$ws.Range("J11").Value = 70004
$ws.Range("K11").Value = 666966.7
$ws.Range("L11").Value = 70004
$ws.Range("M11").Value = -666827.7
$ws.Range("N11").Value = -70282

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9526622
$ws.Range("I16").Value = 3213
$ws.Range("J16").Value = 17859604
$ws.Range("K16").Value = 3213
$ws.Range("L16").Value = 17859604
$ws.Range("M16").Value = -3043
$ws.Range("N16").Value = -17859944
$ws.Range("H22").Value = 9327.615
$ws.Range("J22").Value = 14432.375
$ws.Range("L22").Value = 14432.375
$ws.Range("N22").Value = -15022.375
$ws.Range("H27").Value = 9327.615
$ws.Range("J27").Value = 14432.375
$ws.Range("L27").Value = 14432.375
$ws.Range("N27").Value = -14646.375
$ws.Range("H40").Value = 125004216
$ws.Range("I40").Value = 250002930
$ws.Range("K40").Value = 250002930
$ws.Range("M40").Value = -250002794
$ws.Range("H68").Value = 3558.4
$ws.Range("I68").Value = 2141.8462
$ws.Range("J68").Value = 4395.4546
$ws.Range("K68").Value = 2141.8462
$ws.Range("L68").Value = 4395.4546
$ws.Range("M68").Value = -1392.8462
$ws.Range("N68").Value = -5893.4546
$ws.Range("H71").Value = 3558.4
$ws.Range("I71").Value = 2141.8462
$ws.Range("J71").Value = 4395.4546
$ws.Range("K71").Value = 10709.231
$ws.Range("L71").Value = 21977.273
$ws.Range("M71").Value = -6965.231
$ws.Range("N71").Value = -29465.273
$ws.Range("H132").Value = 4076.303
$ws.Range("I132").Value = 3487.0476
$ws.Range("J132").Value = 5107.5
$ws.Range("K132").Value = 10461.1428
$ws.Range("L132").Value = 15322.5
$ws.Range("M132").Value = -7931.1428
$ws.Range("N132").Value = -20382.5
$ws.Range("H136").Value = 11112835
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 13890544
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 41671632
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -41676732

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 16143.857
$ws.Range("J18").Value = 17501.166
$ws.Range("L18").Value = 17501.166
$ws.Range("N18").Value = -17847.166
$ws.Range("H41").Value = 14000
$ws.Range("J41").Value = 14000
$ws.Range("L41").Value = 14000
$ws.Range("N41").Value = -14780
$ws.Range("H45").Value = 14097.5
$ws.Range("J45").Value = 17626
$ws.Range("L45").Value = 17626
$ws.Range("N45").Value = -18608
$ws.Range("H108").Value = 62600
$ws.Range("J108").Value = 62600
$ws.Range("L108").Value = 62600
$ws.Range("N108").Value = -70280
